# Apply updated benchmark numbers to the comparison table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Simple value updates (style/formatting unchanged) --
$ws.Range("B2").Value  = 0.4348
$ws.Range("B3").Value  = 0.75
$ws.Range("B4").Value  = 0.9592000000000001
$ws.Range("B6").Value  = 0.9314
$ws.Range("B7").Value  = 0.9402
$ws.Range("B8").Value  = 0.9467
$ws.Range("B9").Value  = 0.6381
$ws.Range("I9").Value  = 6
$ws.Range("M9").Value  = 5
$ws.Range("B10").Value = 0.55
$ws.Range("B13").Value = 0.58

# -- Row 5 (Coffee): value + re-highlighted min/max/rank changes --
$ws.Range("B5").Value = 0.8929
$ws.Range("I5").Value = 7
$ws.Range("L5").Value = 5.5
$ws.Range("M5").Value = 5.5

# E5 and F5 are no longer the row minimum, so their "Bad" (red) highlight
# is removed and replaced with the plain/default style used elsewhere in
# the row (e.g. I5). Copy that cell's formatting over instead of inventing
# a new style, so the existing style table entry is reused.
$ws.Range("I5").Copy() | Out-Null
$ws.Range("E5:F5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
